# Time recording log.xlsx edit script
# Commit: "Video 32, 113 teste läksid läbi"
#
# Summary of the change:
#  - A new worksheet "Nädal 8" is appended (a copy of "Nädal 7"'s layout),
#    holding the time entries for the week of 20.03 - 26.03.2020.
#  - "Nädal 7" loses its "currently active" tab state/colour and two of its
#    placeholder "video" markers (rows 9 & 10, col G) are cleared.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate "Nädal 7" (while its tab is still green) to become "Nädal 8".
#    Doing this before recolouring "Nädal 7" means the new sheet naturally
#    inherits the green tab colour that the target file expects.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("Nädal 7")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws8 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8.Name = "Nädal 8"

# ---------------------------------------------------------------------------
# 2) Fix up "Nädal 7": clear the two stray "video" tags left in G9/G10,
#    recolour its tab to pink/red, and move its selection to H11 (it is no
#    longer the active tab once "Nädal 8" is selected below).
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Nädal 7")
$ws7.Range("G9").Value = ""
$ws7.Range("G10").Value = ""
$ws7.Tab.Color = 8420607   # BGR encoding of RGB FFFF7C80
$ws7.Activate()
$ws7.Range("H11").Select()

# ---------------------------------------------------------------------------
# 3) Populate "Nädal 8" with this week's data (week of 20.03 - 26.03.2020).
# ---------------------------------------------------------------------------

# Header: replace the single "week start" date with the explicit date range
# text that was typed in for this particular week.
$ws8.Range("G1").Value = "20.03 - 26.03.2020"

# Row 4 - Friday 20.03.2020, 18:00-19:00, video V32
$ws8.Range("B4").Value = 43910
$ws8.Range("C4").Value = 0.75
$ws8.Range("D4").Value = 0.79166666666666663
$ws8.Range("F4").Value = 60
$ws8.Range("H4").Value = "V32"
$ws8.Range("I4").Value = ""
$ws8.Range("J4").Value = 57

# Row 5 - Saturday 21.03.2020, 14:15-15:30, video V32
$ws8.Range("B5").Value = 43911
$ws8.Range("C5").Value = 0.59375
$ws8.Range("D5").Value = 0.64583333333333337
$ws8.Range("F5").Value = 75
$ws8.Range("H5").Value = "V32"
$ws8.Range("I5").Value = "x"
$ws8.Range("J5").Value = ""

# Row 6 - entry cleared out (kept only the "video" marker)
$ws8.Range("B6").Value = ""
$ws8.Range("C6").Value = ""
$ws8.Range("D6").Value = ""
$ws8.Range("F6").Value = ""
$ws8.Range("H6").Value = ""

# Row 7 - entry cleared out, G7 reset to the generic "video" marker
$ws8.Range("B7").Value = ""
$ws8.Range("C7").Value = ""
$ws8.Range("D7").Value = ""
$ws8.Range("F7").Value = ""
$ws8.Range("G7").Value = "video"
$ws8.Range("H7").Value = ""
$ws8.Range("I7").Value = ""

# Row 8 - entry cleared out
$ws8.Range("C8").Value = ""
$ws8.Range("D8").Value = ""
$ws8.Range("E8").Value = ""
$ws8.Range("F8").Value = ""
$ws8.Range("H8").Value = ""
$ws8.Range("J8").Value = ""

# Row 10 - drop the leftover "video" marker (rows 11+ never had one)
$ws8.Range("G10").Value = ""

# ---------------------------------------------------------------------------
# 4) "Nädal 8" becomes the active/selected tab.
# ---------------------------------------------------------------------------
$ws8.Activate()
$ws8.Range("F6").Select()
